$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("Part 1") changes ---
# Header row: the "Actual" column becomes "Points", and the separate
# "Possible" column of values is no longer tracked (values/formulas cleared,
# but cell styling is left intact).
$ws1.Range("B2").Value = "Points"
$ws1.Range("C2").ClearContents()
$ws1.Range("C3:C11").ClearContents()

$ws1.Range("B14").Value = "Points"
$ws1.Range("C14").ClearContents()
$ws1.Range("C15:C22").ClearContents()

# --- Sheet2 ("Part 2") selection update ---
$ws2.Range("A10").Select() | Out-Null

# --- Sheet1 becomes the active/selected sheet, with its own selection + zoom ---
$ws1.Activate() | Out-Null
$ws1.Range("C24").Select() | Out-Null
$excel.ActiveWindow.Zoom = 150
